$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 (swap with row 9 content)
$ws.Range("D3").Value = 44602
$ws.Range("N3").Value = 6000
$ws.Range("O3").Value = 7000
$ws.Range("P3").Value = 6500
$ws.Range("S3").Value = 3250

# Row 4 (swap with row 10 content)
$ws.Range("D4").Value = 44602
$ws.Range("N4").Value = 5000
$ws.Range("O4").Value = 5000
$ws.Range("P4").Value = 5000
$ws.Range("S4").Value = 2500

# Row 7 (swap with row 13 content)
$ws.Range("D7").Value = 44195
$ws.Range("N7").Value = 3000
$ws.Range("O7").Value = 3500
$ws.Range("P7").Value = 3250
$ws.Range("S7").Value = 1625

# Row 8 (swap with row 14 content)
$ws.Range("D8").Value = 44195
$ws.Range("N8").Value = 2500
$ws.Range("O8").Value = 2500
$ws.Range("P8").Value = 2500
$ws.Range("S8").Value = 1250

# Row 9 (swap with row 3 content)
$ws.Range("D9").Value = 44574
$ws.Range("N9").Value = 7000
$ws.Range("O9").Value = 8000
$ws.Range("P9").Value = 7500
$ws.Range("S9").Value = 3750

# Row 10 (swap with row 4 content)
$ws.Range("D10").Value = 44574
$ws.Range("N10").Value = 6000
$ws.Range("O10").Value = 6000
$ws.Range("P10").Value = 6000
$ws.Range("S10").Value = 3000

# Row 13 (swap with row 7 content)
$ws.Range("D13").Value = 44216
$ws.Range("N13").Value = 3500
$ws.Range("O13").Value = 4000
$ws.Range("P13").Value = 3750
$ws.Range("S13").Value = 1875

# Row 14 (swap with row 8 content)
$ws.Range("D14").Value = 44216
$ws.Range("N14").Value = 3000
$ws.Range("O14").Value = 3000
$ws.Range("P14").Value = 3000
$ws.Range("S14").Value = 1500
